# edit.ps1
# Applies the "gh-pages output generated at 456a3b4" update to 广州-漫展信息.xlsx
#
# Sheet 1 = 展览 (Exhibitions)          -> F-column ("想去人数" / want-to-go count) bumps
# Sheet 2 = 演出 (Performances)          -> F-column bumps
# Sheet 3 = 本地生活 (Local life)        -> untouched
# Sheet 4 = 全部类型 (All categories)    -> F-column bumps + one new row inserted
#           (duplicate "KANAKO ITO&AYANE" entry), shifting the trailing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "展览" - update F column (want-to-go counts)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value  = 494
$ws1.Range("F3").Value  = 1489
$ws1.Range("F4").Value  = 775
$ws1.Range("F5").Value  = 189
$ws1.Range("F6").Value  = 45
$ws1.Range("F7").Value  = 1076
$ws1.Range("F8").Value  = 655
$ws1.Range("F9").Value  = 743
$ws1.Range("F10").Value = 1299
$ws1.Range("F11").Value = 261
$ws1.Range("F12").Value = 998
$ws1.Range("F13").Value = 51
$ws1.Range("F14").Value = 175
$ws1.Range("F16").Value = 389
$ws1.Range("F17").Value = 114
$ws1.Range("F18").Value = 286
$ws1.Range("F19").Value = 522
$ws1.Range("F20").Value = 539
$ws1.Range("F21").Value = 733
$ws1.Range("F23").Value = 152
$ws1.Range("F24").Value = 358

# ---------------------------------------------------------------------------
# Sheet 2: "演出" - update F column (want-to-go counts)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F3").Value  = 969
$ws2.Range("F5").Value  = 215
$ws2.Range("F6").Value  = 13
$ws2.Range("F7").Value  = 129
$ws2.Range("F9").Value  = 565
$ws2.Range("F10").Value = 10

# ---------------------------------------------------------------------------
# Sheet 4: "全部类型" - update F column (want-to-go counts) for the rows that
# existed before the insert (rows 2-34 in the original layout).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value  = 494
$ws4.Range("F4").Value  = 1489
$ws4.Range("F6").Value  = 775
$ws4.Range("F7").Value  = 189
$ws4.Range("F8").Value  = 969
$ws4.Range("F9").Value  = 45
$ws4.Range("F10").Value = 1076
$ws4.Range("F11").Value = 655
$ws4.Range("F12").Value = 743
$ws4.Range("F13").Value = 1299
$ws4.Range("F14").Value = 261
$ws4.Range("F15").Value = 998
$ws4.Range("F16").Value = 51
$ws4.Range("F17").Value = 175
$ws4.Range("F19").Value = 389
$ws4.Range("F21").Value = 215
$ws4.Range("F22").Value = 114
$ws4.Range("F23").Value = 286
$ws4.Range("F24").Value = 13
$ws4.Range("F25").Value = 129
$ws4.Range("F26").Value = 129
$ws4.Range("F27").Value = 522
$ws4.Range("F28").Value = 539
$ws4.Range("F29").Value = 733
$ws4.Range("F32").Value = 152
$ws4.Range("F33").Value = 565
# Row 34 ("KANAKO ITO&AYANE") want-to-go count jumps from 1 to 10 - this is the
# event that gets a freshly-scraped duplicate row inserted right after it.
$ws4.Range("F34").Value = 10

# ---------------------------------------------------------------------------
# Sheet 4: insert a new row 35 - a duplicate "KANAKO ITO&AYANE" entry that sits
# right under row 34, pushing the former rows 35-37 down to 36-38 and growing
# the sheet's dimension from A1:I37 to A1:I38.
# ---------------------------------------------------------------------------
$ws4.Rows.Item(35).Insert()

# Copy cell A34's formatting onto the freshly inserted A35 so the numbering
# column keeps its usual bordered/centered style (s="1") instead of the blank
# default style Excel assigns to a brand new row. (Copy just the single cell,
# not the whole row, so the sheet's used range doesn't balloon out to XFD.)
$ws4.Range("A34").Copy()
$ws4.Range("A35").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws4.Range("A35").Value = 34
# B35 looks like a plain date ("2024-03-31") - force it to be stored as text
# (matching the source data's inlineStr cells) instead of letting Excel's
# autodetection turn it into a date serial number.
$ws4.Range("B35").NumberFormat = "@"
$ws4.Range("B35").Value = "2024-03-31"
$ws4.Range("C35").Value = "【大会员抢先购】广州·KANAKO ITO&AYANE 2024 LIVE"
$ws4.Range("D35").Value = "奥体南路12号优托邦购物中心 疆进酒Omni Space GZ"
$ws4.Range("E35").Value = "2024.03.31 19:00-03.31 20:30"
$ws4.Range("F35").Value = 10
$ws4.Range("G35").Value = 380
$ws4.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=81422"
$ws4.Range("I35").Value = "//i0.hdslb.com/bfs/openplatform/202401/4Y4U8tC01706172039039.jpeg"

# The old row 36 (formerly row 35, "Arknights Only") also got its want-to-go
# count refreshed, from 356 to 358, once it landed in its new position.
$ws4.Range("F37").Value = 358

Write-Host "edit complete"
